# Facility Staffing Combined - "Universal Code edit 2"
#
# 1. Re-labels the quarterly-average column header and adds a new
#    "Staff Increase/Decrease" column (L) computed as (J-B)/J.
# 2. Corrects the state-name labels in column A (proper-case names,
#    "District of Columbia" instead of "D.C.", a couple of labels
#    swapped to the right row) while leaving every row's numeric data
#    untouched.
# 3. Re-enters the K43:K52 average formulas as one block so they are
#    stored as a single shared formula, matching the saved workbook.
# 4. Cosmetic bits: column widths for K/L, and the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header row -----------------------------------------------------
$ws.Range("A1").Value = "State"
$ws.Range("B1").Value = "Q1 2021"
$ws.Range("C1").Value = "Q2 2021"
$ws.Range("D1").Value = "Q3 2021"
$ws.Range("E1").Value = "Q4 2021 "
$ws.Range("F1").Value = "Q1 2022"
$ws.Range("G1").Value = "Q2 2022"
$ws.Range("H1").Value = "Q3 2022"
$ws.Range("I1").Value = "Q4 2022"
$ws.Range("J1").Value = "Q1 2023"
$ws.Range("K1").Value = "Staff Average "
$ws.Range("L1").Value = "Staff Increase/Decrease"

# --- 2. Column A state names --------------------------------------------
$stateNames = @{
    2 = 'Alaska'
    3 = 'Alabama'
    4 = 'Arizona'
    5 = 'Arkansas'
    6 = 'California'
    7 = 'Colorado'
    8 = 'Connecticut'
    9 = 'Delaware'
    10 = 'District of Columbia'
    11 = 'Florida'
    12 = 'Georgia'
    13 = 'Hawaii'
    14 = 'Idaho'
    15 = 'Illinois'
    16 = 'Indiana'
    17 = 'Iowa'
    18 = 'Kansas'
    19 = 'Kentucky'
    20 = 'Louisiana'
    21 = 'Maine'
    22 = 'Maryland'
    23 = 'Massachusetts'
    24 = 'Michigan'
    25 = 'Minnesota'
    26 = 'Mississippi'
    27 = 'Missouri'
    28 = 'Montana'
    29 = 'Nebraska'
    30 = 'Nevada'
    31 = 'New Hampshire'
    32 = 'New Jersey'
    33 = 'New Mexico'
    34 = 'New York'
    35 = 'North Carolina'
    36 = 'North Dakota'
    37 = 'Ohio'
    38 = 'Oklahoma'
    39 = 'Oregon'
    40 = 'Pennsylvania'
    41 = 'Rhode Island'
    42 = 'South Carolina'
    43 = 'South Dakota'
    44 = 'Tennessee'
    45 = 'Texas'
    46 = 'Utah'
    47 = 'Vermont'
    48 = 'Virginia'
    49 = 'Washington'
    50 = 'West Virginia'
    51 = 'Wisconsin'
    52 = 'Wyoming'
}

foreach ($r in $stateNames.Keys) {
    $ws.Range("A$r").Value = $stateNames[$r]
}

# --- 3. New column L: Staff Increase/Decrease ---------------------------
$ws.Range("L2").Formula = "=(J2-B2)/J2"
$ws.Range("L3:L52").Formula = "=(J3-B3)/J3"
$ws.Range("L2:L52").NumberFormat = "0.00"

# --- 4. Re-group the K43:K52 average formulas into one shared formula ---
$ws.Range("K43:K52").Formula = "=(B43+C43+D43+E43+F43+G43+H43+I42+J42)/9"

# --- 5. Cosmetic: column widths & active selection -----------------------
$ws.Columns.Item(11).ColumnWidth = 12.59
$ws.Columns.Item(12).ColumnWidth = 21.59

$ws.Range("C6").Select() | Out-Null
